# PO Status update via batch
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PO_StatusQry")

# Update the "Last updated" timestamp banner in A1
$ws.Range("A1").Value = "Last updated: 2025-07-17 14:43:21"

# Update row 27 (Global PO 4516351202) stats
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 5
$ws.Range("I27").Value = -1
